$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data describing the updated "scheme" endpoint
$ws.Range("A24").Value = "GET"
$ws.Range("B24").Value = "/household/family/list/scheme?type=yolo,student&householdSize=8&totalIncome=300000"
$ws.Range("I24").Value = "Get grant or grants"
$ws.Range("M24").Value = "householdSize, total_income, type={student, family, elder, baby, yolo}"

# Update the view state (scroll position / active selection) as captured in the saved file
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("J9").Select()
